$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add "Area" / "Atotal" columns (G,H) next to the existing x/depth/
# velocity/segment/Q/Qtotal table (A11:F24), plus a small "echo" block
# in J11:K12 that surfaces the new Atotal alongside the existing Qtotal.
# ---------------------------------------------------------------------

# Header row (row 11)
$ws.Range("G11").Value = "Area"
$ws.Range("H11").Value = "Atotal"
$ws.Range("J11").Value = "Atotal"
$ws.Range("K11").Value = "Qtotal"

# Row 12 (summary / first segment row)
$ws.Range("G12").Formula = "=(D12-0)*B12/100"
$ws.Range("H12").Formula = "=SUM(G12:G21)"
$ws.Range("J12").Formula = "=H12"
$ws.Range("K12").Formula = "=F12"

# Row 13 (its own distinct formula, matching the E13 pattern)
$ws.Range("G13").Formula = "=(D13-D12)*B13/100"

# Rows 14:25 share one formula (mirrors the existing D/E shared-formula
# columns), anchored at G14 and filled down through G25
$ws.Range("G14:G25").Formula = "=(D14-D13)*B14/100"

# Selection / active cell ends on the new Atotal/Qtotal echo cells
[void]$ws.Range("J12:K12").Select()
